# Insert 3 new rows before row 301 (shifts old rows 301:385 down to 304:388)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A301:T303").EntireRow.Insert()

# --- New row 301 ---
$ws.Range("A301").Value = 5
$ws.Range("B301").Value = "Macroferia Regional de Talca"
$ws.Range("C301").Value = "Maule"
$ws.Range("D301").Value = 44637
$ws.Range("E301").Value = 7
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100109
$ws.Range("H301").Value = "Uva"
$ws.Range("I301").Value = 100109001
$ws.Range("J301").Value = "Uva"
$ws.Range("K301").Value = "Crimpson Seedless"
$ws.Range("L301").Value = "Primera"
$ws.Range("M301").Value = 300
$ws.Range("N301").Value = 10000
$ws.Range("O301").Value = 10000
$ws.Range("P301").Value = 10000
$ws.Range("Q301").Value = "`$/bandeja 18 kilos"
$ws.Range("R301").Value = "Región de O'Higgins"
$ws.Range("S301").Value = 556
$ws.Range("T301").Value = 18

# --- New row 302 ---
$ws.Range("A302").Value = 5
$ws.Range("B302").Value = "Macroferia Regional de Talca"
$ws.Range("C302").Value = "Maule"
$ws.Range("D302").Value = 44637
$ws.Range("E302").Value = 7
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100109
$ws.Range("H302").Value = "Uva"
$ws.Range("I302").Value = 100109001
$ws.Range("J302").Value = "Uva"
$ws.Range("K302").Value = "Red Globe"
$ws.Range("L302").Value = "Primera"
$ws.Range("M302").Value = 400
$ws.Range("N302").Value = 10000
$ws.Range("O302").Value = 10000
$ws.Range("P302").Value = 10000
$ws.Range("Q302").Value = "`$/bandeja 18 kilos"
$ws.Range("R302").Value = "Provincia de Limarí"
$ws.Range("S302").Value = 556
$ws.Range("T302").Value = 18

# --- New row 303 ---
$ws.Range("A303").Value = 5
$ws.Range("B303").Value = "Macroferia Regional de Talca"
$ws.Range("C303").Value = "Maule"
$ws.Range("D303").Value = 44637
$ws.Range("E303").Value = 7
$ws.Range("F303").Value = "Fruta"
$ws.Range("G303").Value = 100109
$ws.Range("H303").Value = "Uva"
$ws.Range("I303").Value = 100109001
$ws.Range("J303").Value = "Uva"
$ws.Range("K303").Value = "Thompson seedless"
$ws.Range("L303").Value = "Primera"
$ws.Range("M303").Value = 300
$ws.Range("N303").Value = 10000
$ws.Range("O303").Value = 10000
$ws.Range("P303").Value = 10000
$ws.Range("Q303").Value = "`$/bandeja 18 kilos"
$ws.Range("R303").Value = "Región de O'Higgins"
$ws.Range("S303").Value = 556
$ws.Range("T303").Value = 18
